# Renamed con to ctr & added concentration symbol
#
# This script renames every "chb_con2" / "chs_con4" / "chs_con6" /
# "leaf.chb_con2" / "air.chb_con2" column-header label (the "control"
# chamber columns) to the equivalent "ctr" spelling throughout the
# workbook ("chb_ctr2", "chs_ctr4", "chs_ctr6", "leaf.chb_ctr2",
# "air.chb_ctr2"). It also removes the two stray footnote cells on the
# "par" sheet that explained the old "4"/"16" numbering, and finally
# restores the selections / active sheet recorded by Excel at save time.

$wb = $excel.ActiveWorkbook

# --- "chambers" sheet -------------------------------------------------
$chambers = $wb.Worksheets.Item("chambers")
$chambers.Range("E1").Value = "chb_ctr2"
$chambers.Range("G1").Value = "chs_ctr4"
$chambers.Range("I1").Value = "chs_ctr6"

# --- "tc" sheet ---------------------------------------------------------
$tc = $wb.Worksheets.Item("tc")
$tc.Range("D1").Value = "leaf.chb_ctr2"
$tc.Range("E1").Value = "air.chb_ctr2"

# --- "par" sheet ---------------------------------------------------------
$par = $wb.Worksheets.Item("par")
$par.Range("C1").Value = "chb_ctr2"
$par.Range("E1").Value = "chs_ctr4"
$par.Range("G1").Value = "chs_ctr6"
# Drop the old "4 is irrigation" / "16 is control" footnotes entirely.
$par.Range("H4").ClearContents()
$par.Range("H5").ClearContents()

# --- "flow" sheet ---------------------------------------------------------
$flow = $wb.Worksheets.Item("flow")
$flow.Range("C1").Value = "chb_ctr2"
$flow.Range("E1").Value = "chs_ctr4"
$flow.Range("G1").Value = "chs_ctr6"

# --- restore per-sheet selections --------------------------------------
$chambers.Range("I3").Select()
$par.Range("H4:H5").Select()
$flow.Range("G1").Select()

# "tc" becomes the active (visible) tab, so select it last.
$tc.Range("E1").Select()
